$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume table refresh.
# Some 'Price' column values are plain decimal numbers as text (e.g. '2.00', '0.998',
# '1.01'). Assigning such a string straight to .Value lets Excel's input parser treat
# it as a number and silently renormalize it (dropping significant trailing zeros,
# switching to scientific notation, etc.), which would corrupt the source formatting.
# Forcing NumberFormat to Text ('@') first keeps those cells exactly as authored.

# Row 2
$ws.Range("D2").Value = '67.535.53'
$ws.Range("E2").Value = '  -3.10%  '

# Row 3
$ws.Range("D3").Value = '3.797.62'
$ws.Range("E3").Value = '  +1.48%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.57'
$ws.Range("E5").Value = '  -3.41%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.99'
$ws.Range("E6").Value = '  -3.83%  '

# Row 7
$ws.Range("D7").Value = '3.787.11'
$ws.Range("E7").Value = '  +1.18%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").Value = '  -1.45%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'
$ws.Range("E10").Value = '  -4.18%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.22'
$ws.Range("E11").Value = '  -5.26%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  -4.15%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.74'
$ws.Range("E13").Value = '  -5.48%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'
$ws.Range("E14").Value = '  -3.86%  '

# Row 15
$ws.Range("D15").Value = '4.416.20'
$ws.Range("E15").Value = '  +1.32%  '

# Row 16
$ws.Range("D16").Value = '3.777.66'
$ws.Range("E16").Value = '  +1.04%  '

# Row 17
$ws.Range("D17").Value = '67.574.44'

# Row 18
$ws.Range("E18").Value = '  -4.77%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.14'
$ws.Range("E19").Value = '  -4.44%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.04'
$ws.Range("E20").Value = '  -1.90%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '487.62'
$ws.Range("E21").Value = '  -2.89%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.19'
$ws.Range("E22").Value = '  +0.55%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("E23").Value = '  +0.19%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.94'
$ws.Range("E24").Value = '  -2.49%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  -11.10%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000141'
$ws.Range("E26").Value = '  +3.28%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.18'
$ws.Range("E27").Value = '  -5.62%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.18'
$ws.Range("E28").Value = '  -12.44%  '

# Row 29
$ws.Range("E29").Value = '  -0.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.91'
$ws.Range("E30").Value = '  -0.51%  '

# Row 31
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.68'
$ws.Range("E31").Value = '  +7.54%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.40'
$ws.Range("E32").Value = '  -3.27%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.79'
$ws.Range("E33").Value = '  -3.15%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.108'
$ws.Range("E34").Value = '  -4.87%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  -0.18%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("E36").Value = '  -3.46%  '

# Row 37
$ws.Range("E37").Value = '  -1.36%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.73'
$ws.Range("E38").Value = '  -6.29%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '451.17'
$ws.Range("E39").Value = '  +1.76%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.323'
$ws.Range("E40").Value = '  -9.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '48.87'
$ws.Range("E41").Value = '  -1.72%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.00'
$ws.Range("E42").Value = '  -3.66%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.85'
$ws.Range("E43").Value = '  -7.04%  '

# Row 44
$ws.Range("B44").Value = 'Cosmos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.26'
$ws.Range("E44").Value = '  -3.53%  '

# Row 45
$ws.Range("B45").Value = 'Arweave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.47'
$ws.Range("E45").Value = '  -10.13%  '

# Row 46
$ws.Range("D46").Value = '2.835.94'
$ws.Range("E46").Value = '  -4.02%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.56'
$ws.Range("E47").Value = '  +1.39%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0349'
$ws.Range("E49").Value = '  -3.06%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.94'
$ws.Range("E50").Value = '  -4.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.42'
$ws.Range("E51").Value = '  +8.10%  '
